$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ----------------------------------------------------------------------
# Change 1: close the inner "IF A > C" block with its own "END IF" line,
# inserted right before the "	ELSE " paragraph that closes the outer
# "IF A > B" block.
# ----------------------------------------------------------------------
$pElse = $d.Paragraphs.Item(12)
[void]$pElse.Range.InsertParagraphBefore()
$pNewEndIf = $d.Paragraphs.Item(12)
[void]$pNewEndIf.Range.InsertXML("<w:p $wNs><w:r><w:tab/></w:r><w:r><w:tab/><w:t>END IF</w:t></w:r></w:p>")

# ----------------------------------------------------------------------
# Change 2: the "END IF" line that closes the "IF B > C" block gets a
# second (nesting) tab, matching the indentation of the other nested
# blocks.
# ----------------------------------------------------------------------
$pEndIf = $d.Paragraphs.Item(18)
[void]$pEndIf.Range.InsertXML("<w:p $wNs><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:t>END IF</w:t></w:r></w:p>")

# ----------------------------------------------------------------------
# Change 3: the final "END" paragraph becomes "	END IF" (keeping the
# _GoBack bookmark in place), and change 4 adds a brand-new paragraph
# right after it containing just "END".
# ----------------------------------------------------------------------
$pEnd = $d.Paragraphs.Item(19)
[void]$pEnd.Range.InsertXML("<w:p $wNs><w:r><w:tab/><w:t>END IF</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>")

$pEndIf2 = $d.Paragraphs.Item(19)
[void]$pEndIf2.Range.InsertParagraphAfter()
$pNewEnd = $d.Paragraphs.Item(20)
[void]$pNewEnd.Range.InsertXML("<w:p $wNs><w:r><w:t>END</w:t></w:r></w:p>")
